# Fill in the newly-computed proposed_com / proposed_emb results (columns I/J)
# for the Random Forest, Naive Bayes and Neural Network blocks. These cells
# previously held the placeholder text "-"; once every "-" cell is replaced
# the shared string is no longer referenced anywhere in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.0234744184665354
$ws.Range("J2").Value = 0.0135729788102669
$ws.Range("I3").Value = 0.0259999999999999
$ws.Range("J3").Value = 0.026
$ws.Range("I4").Value = 0.0212069180508816
$ws.Range("J4").Value = 0.0243603282906484
$ws.Range("I5").Value = 0.032
$ws.Range("J5").Value = 0.034
$ws.Range("I6").Value = 0.0278223478109327
$ws.Range("J6").Value = 0.0305382356169273
$ws.Range("I7").Value = 0.0424
$ws.Range("J7").Value = 0.0396
$ws.Range("I8").Value = 0.0374271052252067
$ws.Range("J8").Value = 0.0361822593958852
$ws.Range("I9").Value = 0.00327450076427656
$ws.Range("J9").Value = 0.00195748351218722
$ws.Range("I18").Value = 0.0381760904894651
$ws.Range("J18").Value = 0.0331001953043007
$ws.Range("I19").Value = 0.054
$ws.Range("J19").Value = 0.04
$ws.Range("I20").Value = 0.0591263286310486
$ws.Range("J20").Value = 0.0417276388330926
$ws.Range("I21").Value = 0.054
$ws.Range("J21").Value = 0.044
$ws.Range("I22").Value = 0.0570345387681469
$ws.Range("J22").Value = 0.04385660329967
$ws.Range("I23").Value = 0.0536
$ws.Range("J23").Value = 0.0416
$ws.Range("I24").Value = 0.0551191478409783
$ws.Range("J24").Value = 0.0421365397678388
$ws.Range("I25").Value = 0.302071880211738
$ws.Range("J25").Value = 0.300734654614856
$ws.Range("I34").Value = 0.0285714285714285
$ws.Range("J34").Value = 0.04
$ws.Range("I35").Value = 0.16
$ws.Range("J35").Value = 0.223999999999999
$ws.Range("I36").Value = 0.164650083870541
$ws.Range("J36").Value = 0.228936671568393
$ws.Range("I37").Value = 0.158
$ws.Range("J37").Value = 0.184
$ws.Range("I38").Value = 0.162660913424754
$ws.Range("J38").Value = 0.19874511529818
$ws.Range("I39").Value = 0.1616
$ws.Range("J39").Value = 0.1472
$ws.Range("I40").Value = 0.163641566110344
$ws.Range("J40").Value = 0.164481635046492
$ws.Range("I41").Value = 0.000765234765234765
$ws.Range("J41").Value = 0.000518047608956699

# Leave the cursor where the author last left it when saving.
$null = $ws.Range("H11").Select()
